$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 264 is an exact duplicate of row 263 (same item.name/url/subclass/
# division/retailer) -- "corrected a duplicated product". Delete the
# duplicated row; every row below it shifts up by one automatically.
$ws.Rows.Item(264).Delete()

# Reflect the author's final cursor position after the cleanup.
$ws.Range("A264").Select()
